# Updated symbol list on Fri Feb 10 04:51:41 UTC 2023 with GitHub Actions
# Re-applies the refreshed coinranking.com price/volume snapshot onto the
# existing Sheet1 table (cells are stored as text in the source data, so we
# force a Text number format before writing, then restore the default style
# so we do not leave a stray numeric/percent format behind).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, in diff order.
$edits = @(
    @{ Ref = 'D2'; Value = '306.82' },
    @{ Ref = 'E2'; Value = '-4.01%' },
    @{ Ref = 'D3'; Value = '39.90' },
    @{ Ref = 'E3'; Value = '-6.65%' },
    @{ Ref = 'D4'; Value = '5.041' },
    @{ Ref = 'E4'; Value = '-3.43%' },
    @{ Ref = 'D5'; Value = '0.07672' },
    @{ Ref = 'E5'; Value = '-5.96%' },
    @{ Ref = 'D6'; Value = '4.233' },
    @{ Ref = 'E6'; Value = '-2.02%' },
    @{ Ref = 'D7'; Value = '1.614' },
    @{ Ref = 'E7'; Value = '-8.57%' },
    @{ Ref = 'D8'; Value = '0.8891' },
    @{ Ref = 'E8'; Value = '-6.57%' },
    @{ Ref = 'D9'; Value = '0.1006' },
    @{ Ref = 'E9'; Value = '-9.52%' },
    @{ Ref = 'D10'; Value = '0.1726' },
    @{ Ref = 'E10'; Value = '-7.03%' },
    @{ Ref = 'D11'; Value = '0.08882' },
    @{ Ref = 'E11'; Value = '-5.37%' },
    @{ Ref = 'D12'; Value = '0.04387' },
    @{ Ref = 'E12'; Value = '-5.89%' },
    @{ Ref = 'E13'; Value = '-0.32%' },
    @{ Ref = 'B14'; Value = 'TigerCash' },
    @{ Ref = 'C14'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' },
    @{ Ref = 'D14'; Value = '0.005805' },
    @{ Ref = 'E14'; Value = '-0.43%' },
    @{ Ref = 'B15'; Value = 'LEO' },
    @{ Ref = 'C15'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' },
    @{ Ref = 'D15'; Value = '3.355' },
    @{ Ref = 'E15'; Value = '-0.40%' },
    @{ Ref = 'B16'; Value = 'BTSEToken' },
    @{ Ref = 'C16'; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' },
    @{ Ref = 'D16'; Value = '2.530' },
    @{ Ref = 'E16'; Value = '0.53%' },
    @{ Ref = 'B17'; Value = 'BitpandaEcosystemToken' },
    @{ Ref = 'C17'; Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best' },
    @{ Ref = 'D17'; Value = '0.3361' },
    @{ Ref = 'E17'; Value = '-0.06%' },
    @{ Ref = 'B18'; Value = 'MCDex' },
    @{ Ref = 'C18'; Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb' },
    @{ Ref = 'D18'; Value = '7.089' },
    @{ Ref = 'E18'; Value = '-4.80%' },
    @{ Ref = 'B19'; Value = 'ProBitToken' },
    @{ Ref = 'C19'; Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob' },
    @{ Ref = 'D19'; Value = '0.1342' },
    @{ Ref = 'E19'; Value = '-3.36%' },
    @{ Ref = 'B20'; Value = 'ZBToken' },
    @{ Ref = 'C20'; Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb' },
    @{ Ref = 'D20'; Value = '0.3207' },
    @{ Ref = 'E20'; Value = '19.37%' },
    @{ Ref = 'B21'; Value = 'BitForexToken' },
    @{ Ref = 'C21'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' },
    @{ Ref = 'D21'; Value = '0.001269' },
    @{ Ref = 'E21'; Value = '-0.68%' },
    @{ Ref = 'D22'; Value = '0.04217' },
    @{ Ref = 'E22'; Value = '0.87%' },
    @{ Ref = 'D23'; Value = '0.001195' },
    @{ Ref = 'E23'; Value = '-4.46%' },
    @{ Ref = 'D24'; Value = '0.004061' },
    @{ Ref = 'E24'; Value = '-5.94%' },
    @{ Ref = 'E25'; Value = '-6.49%' },
    @{ Ref = 'E26'; Value = '-0.01%' },
    @{ Ref = 'D38'; Value = '0.02344' },
    @{ Ref = 'E38'; Value = '-9.75%' },
    @{ Ref = 'D39'; Value = '0.05156' },
    @{ Ref = 'E39'; Value = '-6.04%' },
    @{ Ref = 'D40'; Value = '0.007949' },
    @{ Ref = 'E40'; Value = '2.22%' },
    @{ Ref = 'D41'; Value = '0.1321' },
    @{ Ref = 'E41'; Value = '-5.15%' },
    @{ Ref = 'D42'; Value = '0.006559' },
    @{ Ref = 'E42'; Value = '-0.74%' },
    @{ Ref = 'E43'; Value = '-6.61%' },
    @{ Ref = 'D44'; Value = '0.008670' },
    @{ Ref = 'E44'; Value = '2.18%' },
    @{ Ref = 'D45'; Value = '0.3043' },
    @{ Ref = 'E45'; Value = '-11.87%' },
    @{ Ref = 'D46'; Value = '0.00006569' },
    @{ Ref = 'E46'; Value = '-5.67%' },
    @{ Ref = 'E47'; Value = '0.08%' },
    @{ Ref = 'D48'; Value = '0.003404' },
    @{ Ref = 'E48'; Value = '-2.12%' },
    @{ Ref = 'E49'; Value = '41.58%' },
    @{ Ref = 'D50'; Value = '0.00002104' },
    @{ Ref = 'E50'; Value = '0.08%' },
    @{ Ref = 'D51'; Value = '0.0002004' },
    @{ Ref = 'E51'; Value = '0.08%' }
)

foreach ($edit in $edits) {
    $col = $edit.Ref -replace '[0-9]+$', ''
    $range = $ws.Range($edit.Ref)
    if ($col -eq "D" -or $col -eq "E") {
        # Price / Volume(1h) columns hold numeric- and percent-looking text;
        # force Text format so Excel does not reinterpret the string as a number.
        $range.NumberFormat = "@"
        $range.Value = $edit.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $edit.Value
    }
}
